# Cowboy Treasure Deluxe review: move the "Meta description" blurb from
# right under the title down to the very end of the document, where it
# replaces the old "Prompt: ..." paragraph. The title line is restated
# (bold) as a lead-in just above it.

$d = $word.ActiveDocument

# Step 1: Remove the "Meta description" paragraph (currently the 2nd
# paragraph, right after the H1 title).
$null = $d.Paragraphs(2).Range.Delete()

# Step 2: Replace the final "Prompt: ..." paragraph with two new
# paragraphs: a bold restatement of the page title, followed by an
# italic paragraph holding the description text that used to follow
# "Meta description: " near the top.
$lastP = $d.Paragraphs($d.Paragraphs.Count)
$targetRange = $d.Range($lastP.Range.Start, $lastP.Range.End - 1)

$xmlSnippet = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Cowboy Treasure Deluxe Free - Pros and Cons | Review</w:t></w:r></w:p><w:p><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Explore the Wild West and find out the pros and cons of Cowboy Treasure Deluxe online slot. Play now for free.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$null = $targetRange.InsertXML($xmlSnippet)
